$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Team")

# Update column G (totaltimetaken) values from 25 to 20 for rows 2-11
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 7).Value = 20
}

# Update the selected cell/range shown in the sheet view to G12
$ws.Activate()
$ws.Range("G12").Select()
